$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 94.45181500000001
$ws.Range("H2").Value = 283.355445
$ws.Range("I2").Value = 0.1457155743604623
$ws.Range("J2").Value = 0.1590548236774281
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.495057333333333
$ws.Range("N2").Value = 16.485172
$ws.Range("O2").Value = 0.8161989011161211
$ws.Range("P2").Value = 0.8403205285996808
$ws.Range("Q2").Value = 519.0181386623933
$ws.Range("R2").Value = 4671.16324796154
$ws.Range("S2").Value = 0.1189328916685138
$ws.Range("T2").Value = 0.1336570335089454

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 94.45181500000001
$ws.Range("H3").Value = 283.355445
$ws.Range("I3").Value = 0.1457155743604623
$ws.Range("J3").Value = 0.1590548236774281
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.657666
$ws.Range("N3").Value = 1.972998
$ws.Range("O3").Value = 0.09768528951377062
$ws.Range("P3").Value = 0.1005722428790014
$ws.Range("Q3").Value = 62.11774736379
$ws.Range("R3").Value = 559.0597262741101
$ws.Range("S3").Value = 0.01423426806806713
$ws.Range("T3").Value = 0.01599650035796304

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 94.45181500000001
$ws.Range("H4").Value = 283.355445
$ws.Range("I4").Value = 0.1457155743604623
$ws.Range("J4").Value = 0.1590548236774281
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 0.5797745
$ws.Range("N4").Value = 1.159549
$ws.Range("O4").Value = 0.08611580937010824
$ws.Range("P4").Value = 0.0591072285213179
$ws.Range("Q4").Value = 54.7607538157175
$ws.Range("R4").Value = 328.564522894305
$ws.Range("S4").Value = 0.01254841462388141
$ws.Range("T4").Value = 0.009401289810519667

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 125.4169656666667
$ws.Range("H5").Value = 376.250897
$ws.Range("I5").Value = 0.1934870726059072
$ws.Range("J5").Value = 0.2111994709712006
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.495057333333333
$ws.Range("N5").Value = 16.485172
$ws.Range("O5").Value = 0.8161989011161211
$ws.Range("P5").Value = 0.8403205285996808
$ws.Range("Q5").Value = 689.1734169110315
$ws.Range("R5").Value = 6202.560752199283
$ws.Range("S5").Value = 0.1579239360411166
$ws.Range("T5").Value = 0.1774752510864923

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 125.4169656666667
$ws.Range("H6").Value = 376.250897
$ws.Range("I6").Value = 0.1934870726059072
$ws.Range("J6").Value = 0.2111994709712006
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.657666
$ws.Range("N6").Value = 1.972998
$ws.Range("O6").Value = 0.09768528951377062
$ws.Range("P6").Value = 0.1005722428790014
$ws.Range("Q6").Value = 82.482474142134
$ws.Range("R6").Value = 742.342267279206
$ws.Range("S6").Value = 0.01890084070468
$ws.Range("T6").Value = 0.0212408044904322

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 125.4169656666667
$ws.Range("H7").Value = 376.250897
$ws.Range("I7").Value = 0.1934870726059072
$ws.Range("J7").Value = 0.2111994709712006
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 0.5797745
$ws.Range("N7").Value = 1.159549
$ws.Range("O7").Value = 0.08611580937010824
$ws.Range("P7").Value = 0.0591072285213179
$ws.Range("Q7").Value = 72.71355856090884
$ws.Range("R7").Value = 436.281351365453
$ws.Range("S7").Value = 0.0166622958601106
$ws.Range("T7").Value = 0.0124834153942762

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 129.124439
$ws.Range("H8").Value = 387.373317
$ws.Range("I8").Value = 0.1992067785341921
$ws.Range("J8").Value = 0.2174427762726615
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.495057333333333
$ws.Range("N8").Value = 16.485172
$ws.Range("O8").Value = 0.8161989011161211
$ws.Range("P8").Value = 0.8403205285996808
$ws.Range("Q8").Value = 709.5461954395028
$ws.Range("R8").Value = 6385.915758955524
$ws.Range("S8").Value = 0.1625923537344901
$ws.Range("T8").Value = 0.1827216286976251

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 129.124439
$ws.Range("H9").Value = 387.373317
$ws.Range("I9").Value = 0.1992067785341921
$ws.Range("J9").Value = 0.2174427762726615
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.657666
$ws.Range("N9").Value = 1.972998
$ws.Range("O9").Value = 0.09768528951377062
$ws.Range("P9").Value = 0.1005722428790014
$ws.Range("Q9").Value = 84.92075329937401
$ws.Range("R9").Value = 764.2867796943661
$ws.Range("S9").Value = 0.01945957183421814
$ws.Range("T9").Value = 0.02186870770757848

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 129.124439
$ws.Range("H10").Value = 387.373317
$ws.Range("I10").Value = 0.1992067785341921
$ws.Range("J10").Value = 0.2174427762726615
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 0.5797745
$ws.Range("N10").Value = 1.159549
$ws.Range("O10").Value = 0.08611580937010824
$ws.Range("P10").Value = 0.0591072285213179
$ws.Range("Q10").Value = 74.86305705900551
$ws.Range("R10").Value = 449.178342354033
$ws.Range("S10").Value = 0.01715485296548386
$ws.Range("T10").Value = 0.01285243986745801

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 136.1162363333333
$ws.Range("H11").Value = 408.348709
$ws.Range("I11").Value = 0.2099933765920337
$ws.Range("J11").Value = 0.2292168125052277
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.495057333333333
$ws.Range("N11").Value = 16.485172
$ws.Range("O11").Value = 0.8161989011161211
$ws.Range("P11").Value = 0.8403205285996808
$ws.Range("Q11").Value = 747.9665226492162
$ws.Range("R11").Value = 6731.698703842947
$ws.Range("S11").Value = 0.1713963632160817
$ws.Range("T11").Value = 0.1926155930483269

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 136.1162363333333
$ws.Range("H12").Value = 408.348709
$ws.Range("I12").Value = 0.2099933765920337
$ws.Range("J12").Value = 0.2292168125052277
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.657666
$ws.Range("N12").Value = 1.972998
$ws.Range("O12").Value = 0.09768528951377062
$ws.Range("P12").Value = 0.1005722428790014
$ws.Range("Q12").Value = 89.51902068439799
$ws.Range("R12").Value = 805.671186159582
$ws.Range("S12").Value = 0.02051326378836707
$ws.Range("T12").Value = 0.02305284893922629

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 136.1162363333333
$ws.Range("H13").Value = 408.348709
$ws.Range("I13").Value = 0.2099933765920337
$ws.Range("J13").Value = 0.2292168125052277
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 0.5797745
$ws.Range("N13").Value = 1.159549
$ws.Range("O13").Value = 0.08611580937010824
$ws.Range("P13").Value = 0.0591072285213179
$ws.Range("Q13").Value = 78.91672286204016
$ws.Range("R13").Value = 473.500337172241
$ws.Range("S13").Value = 0.01808374958758492
$ws.Range("T13").Value = 0.01354837051767458

$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 163.083542
$ws.Range("H14").Value = 326.167084
$ws.Range("I14").Value = 0.2515971979074048
$ws.Range("J14").Value = 0.183086116573482
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 5.495057333333333
$ws.Range("N14").Value = 16.485172
$ws.Range("O14").Value = 0.8161989011161211
$ws.Range("P14").Value = 0.8403205285996808
$ws.Range("Q14").Value = 896.1534134130746
$ws.Range("R14").Value = 5376.920480478449
$ws.Range("S14").Value = 0.205353356455919
$ws.Range("T14").Value = 0.1538510222582912

$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 163.083542
$ws.Range("H15").Value = 326.167084
$ws.Range("I15").Value = 0.2515971979074048
$ws.Range("J15").Value = 0.183086116573482
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.657666
$ws.Range("N15").Value = 1.972998
$ws.Range("O15").Value = 0.09768528951377062
$ws.Range("P15").Value = 0.1005722428790014
$ws.Range("Q15").Value = 107.254500732972
$ws.Range("R15").Value = 643.5270043978321
$ws.Range("S15").Value = 0.02457734511843828
$ws.Range("T15").Value = 0.01841338138380139

$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 163.083542
$ws.Range("H16").Value = 326.167084
$ws.Range("I16").Value = 0.2515971979074048
$ws.Range("J16").Value = 0.183086116573482
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 0.5797745
$ws.Range("N16").Value = 1.159549
$ws.Range("O16").Value = 0.08611580937010824
$ws.Range("P16").Value = 0.0591072285213179
$ws.Range("Q16").Value = 94.55167902127901
$ws.Range("R16").Value = 378.2067160851161
$ws.Range("S16").Value = 0.02166649633304746
$ws.Range("T16").Value = 0.01082171293138945
